# Internal version 2.0. TFS 3416
# Replace jobcode WISY13 with WISY14, and log the change in Revision_History.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. AT_Module_Access: rows for jobcode WISY13 ('Sr Analyst, Systems')
#    -> update JobCode column (A) and the documentation/comment column (I)
#    which contains literal SQL VALUES text referencing the jobcode.
# ---------------------------------------------------------------------------
$wsModule = $wb.Worksheets.Item("AT_Module_Access")

$wsModule.Range("A2").Value = "WISY14"
$wsModule.Range("A3").Value = "WISY14"
$wsModule.Range("A4").Value = "WISY14"
$wsModule.Range("A5").Value = "WISY14"
$wsModule.Range("A6").Value = "WISY14"

$wsModule.Range("I8").Value  = "           ('WISY14','Sr Analyst, Systems',1,'CSR',1),"
$wsModule.Range("I9").Value  = "           ('WISY14','Sr Analyst, Systems',2,'Supervisor',1),"
$wsModule.Range("I10").Value = "           ('WISY14','Sr Analyst, Systems',3,'Quality',1),"
$wsModule.Range("I11").Value = "           ('WISY14','Sr Analyst, Systems',4,'LSA',1),"
$wsModule.Range("I12").Value = "           ('WISY14','Sr Analyst, Systems',5,'Training',1),"

$wsModule.Activate()
$wsModule.Range("I12").Select()

# ---------------------------------------------------------------------------
# 2. AT_Role_Access: rows for jobcode WISY13 ('Sr Analyst, Systems')
#    -> update JobCode column (A) and the documentation/comment column (K).
# ---------------------------------------------------------------------------
$wsRole = $wb.Worksheets.Item("AT_Role_Access")

$wsRole.Range("A2").Value = "WISY14"
$wsRole.Range("A8").Value = "WISY14"

$wsRole.Range("K9").Value  = "           ('WISY14','Sr Analyst, Systems',101,'CoachingAdmin',0,1),"
$wsRole.Range("K15").Value = "           ('WISY14','Sr Analyst, Systems',103,'WarningAdmin',0,1)"

$wsRole.Activate()
$wsRole.Range("K15").Select()

# ---------------------------------------------------------------------------
# 3. Revision_History: log this change as revision 2.
# ---------------------------------------------------------------------------
$wsRev = $wb.Worksheets.Item("Revision_History")

$wsRev.Range("A4").Value = 2

# Copy B2's date format onto B4 before assigning the serial date value so the
# new row keeps the same short-date style already used by the other entries.
$wsRev.Range("B2").Copy()
$wsRev.Range("B4").PasteSpecial(-4122)
$wsRev.Range("B4").Value = 42577

$wsRev.Range("C4").Value = "Susmitha Palacherla"
$wsRev.Range("D4").Value = 3416
$wsRev.Range("E4").Value = "Replace jobcode WISY13 with WISY14"

# Leave the workbook's active tab/selection back on Revision_History, as it
# was originally (that sheet stayed the selected tab in the real edit).
$wsRev.Activate()
$wsRev.Range("A1").Select()
